$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new row (82) with the next month's data point: "01-09-2021" / 99.90000000000001
# The date-like text must land as a plain text/shared-string value (matching the
# existing A2:A81 cells), not get auto-converted to a date serial number. Using
# Formula with a quoted text literal keeps it a string; copy+paste-special values
# then collapses the formula down to a literal cached string without touching
# cell styles (no NumberFormat round-trip, so no new style entries are created).
$ws.Range("A82").Formula = '="01-09-2021"'
$ws.Range("A82").Copy() | Out-Null
$ws.Range("A82").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B82").Value = 99.90000000000001
